$d = $word.ActiveDocument

# 1) Update the letter date: "September 19, 2025" -> "September 21, 2025".
#    Set the run's Range.Text directly (rather than Find & Replace) so the
#    existing run/formatting is reused as-is.
$dateRange = $d.Paragraphs(4).Range
if ($dateRange.Text.TrimEnd([char]13) -eq "September 19, 2025") {
    $dateRange.Text = "September 21, 2025"
}

# 2) Split the mailing-address paragraph ("2973 Lamory Pl, Santa Clara CA 95051")
#    into two paragraphs: "2973 Lamory Pl" and "Santa Clara, CA 95051".
#    This exact text also appears later in the "PROPERTY ADDRESS:" table cell,
#    which must stay untouched, so edit paragraph 7 (the mailing address,
#    right under "Moliang Zhou") directly instead of a document-wide replace.
$addressPara = $d.Paragraphs(7)
$addressRange = $addressPara.Range
if ($addressRange.Text.TrimEnd([char]13) -eq "2973 Lamory Pl, Santa Clara CA 95051") {
    $addressRange.Text = "2973 Lamory Pl" + [char]13 + "Santa Clara, CA 95051"
}

# 3) Remove the now-superfluous blank "No Spacing" paragraph that sat directly
#    below "... Board of Directors" (it shifted from index 41 to 42 because of
#    the paragraph inserted by step 2 above).
$blankPara = $d.Paragraphs(42)
if ($blankPara.Range.Text.TrimEnd([char]13) -eq "" -and $blankPara.Style.NameLocal -eq "No Spacing") {
    $blankPara.Range.Delete() | Out-Null
}
